# Auto-generated edit script: updates crypto price/volume table cells
# per the commit "Updated cryptos list on Mon May 20 04:46:07 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.003.79"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "3.112.66"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'578.93"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'174.05"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.112.73"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'6.40"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "'0.476"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "'36.08"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "3.627.54"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "66.897.17"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.02"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.06"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").Value = "3.106.37"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'487.21"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "'7.85"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").Value = "'0.694"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Value = "'83.65"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").Value = "'12.83"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").Value = "'2.29"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").Value = "'2.60"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("D32").Value = "'28.19"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "0.0₃0941"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'48.15"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "'5.63"
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.311"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'49.22"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("D43").Value = "'8.32"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "2.805.09"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "'0.0348"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "'371.94"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").Value = "'134.62"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D50").Value = "'24.57"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'2.22"
$ws.Range("E51").Value = "  +0.20%  "
